$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Add a new row of data (row 4) to the sheet.
$ws.Range("A4").Value = 5959595959599
$ws.Range("B4").Value = "ss"
$ws.Range("C4").Value = "dd"
$ws.Range("D4").Value = "CD"
$ws.Range("E4").Value = 2

# "Release Date" for this row is entered as plain text, not a date value.
$ws.Range("F4").Value = "23.02.2020"

$ws.Range("G4").Value = "DVD02"
